$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44203
$ws.Cells.Item(2, 10).Value = 300
$ws.Cells.Item(2, 11).Value = 4500
$ws.Cells.Item(2, 12).Value = 5000
$ws.Cells.Item(2, 13).Value = 4750
$ws.Cells.Item(2, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(2, 16).Value = 79
$ws.Cells.Item(2, 17).Value = 60

$ws.Cells.Item(3, 4).Value = 44251
$ws.Cells.Item(3, 10).Value = 700
$ws.Cells.Item(3, 11).Value = 6500
$ws.Cells.Item(3, 12).Value = 7000
$ws.Cells.Item(3, 13).Value = 6750
$ws.Cells.Item(3, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(3, 16).Value = 112
$ws.Cells.Item(3, 17).Value = 60

$ws.Cells.Item(4, 4).Value = 44217
$ws.Cells.Item(4, 10).Value = 700
$ws.Cells.Item(4, 11).Value = 6500
$ws.Cells.Item(4, 12).Value = 7000
$ws.Cells.Item(4, 13).Value = 6750
$ws.Cells.Item(4, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(4, 16).Value = 112
$ws.Cells.Item(4, 17).Value = 60

$ws.Cells.Item(5, 4).Value = 44504
$ws.Cells.Item(5, 10).Value = 700
$ws.Cells.Item(5, 11).Value = 6500
$ws.Cells.Item(5, 12).Value = 7000
$ws.Cells.Item(5, 13).Value = 6750
$ws.Cells.Item(5, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(5, 16).Value = 112
$ws.Cells.Item(5, 17).Value = 60

$ws.Cells.Item(6, 4).Value = 44301
$ws.Cells.Item(6, 10).Value = 300
$ws.Cells.Item(6, 11).Value = 6000
$ws.Cells.Item(6, 12).Value = 7000
$ws.Cells.Item(6, 13).Value = 6500
$ws.Cells.Item(6, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(6, 16).Value = 108
$ws.Cells.Item(6, 17).Value = 60

$ws.Cells.Item(7, 4).Value = 44279
$ws.Cells.Item(7, 10).Value = 500
$ws.Cells.Item(7, 11).Value = 7000
$ws.Cells.Item(7, 12).Value = 8000
$ws.Cells.Item(7, 13).Value = 7500
$ws.Cells.Item(7, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(7, 16).Value = 150
$ws.Cells.Item(7, 17).Value = 50

$ws.Cells.Item(8, 4).Value = 44294
$ws.Cells.Item(8, 10).Value = 500
$ws.Cells.Item(8, 11).Value = 7000
$ws.Cells.Item(8, 12).Value = 8000
$ws.Cells.Item(8, 13).Value = 7500
$ws.Cells.Item(8, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(8, 16).Value = 125
$ws.Cells.Item(8, 17).Value = 60

$ws.Cells.Item(9, 4).Value = 44238
$ws.Cells.Item(9, 10).Value = 400
$ws.Cells.Item(9, 11).Value = 7000
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 7500
$ws.Cells.Item(9, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(9, 16).Value = 125
$ws.Cells.Item(9, 17).Value = 60

$ws.Cells.Item(10, 4).Value = 44266
$ws.Cells.Item(10, 10).Value = 600
$ws.Cells.Item(10, 11).Value = 6500
$ws.Cells.Item(10, 12).Value = 7000
$ws.Cells.Item(10, 13).Value = 6750
$ws.Cells.Item(10, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(10, 16).Value = 112
$ws.Cells.Item(10, 17).Value = 60

$ws.Cells.Item(11, 4).Value = 44482
$ws.Cells.Item(11, 10).Value = 400
$ws.Cells.Item(11, 11).Value = 11000
$ws.Cells.Item(11, 12).Value = 12000
$ws.Cells.Item(11, 13).Value = 11500
$ws.Cells.Item(11, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(11, 16).Value = 192
$ws.Cells.Item(11, 17).Value = 60

$ws.Cells.Item(12, 4).Value = 44483
$ws.Cells.Item(12, 10).Value = 340
$ws.Cells.Item(12, 11).Value = 10000
$ws.Cells.Item(12, 12).Value = 11000
$ws.Cells.Item(12, 13).Value = 10500
$ws.Cells.Item(12, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(12, 16).Value = 175
$ws.Cells.Item(12, 17).Value = 60

$ws.Cells.Item(13, 4).Value = 44524
$ws.Cells.Item(13, 10).Value = 508
$ws.Cells.Item(13, 11).Value = 5000
$ws.Cells.Item(13, 12).Value = 6000
$ws.Cells.Item(13, 13).Value = 5508
$ws.Cells.Item(13, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(13, 16).Value = 92
$ws.Cells.Item(13, 17).Value = 60

$ws.Cells.Item(14, 4).Value = 44258
$ws.Cells.Item(14, 10).Value = 500
$ws.Cells.Item(14, 11).Value = 7000
$ws.Cells.Item(14, 12).Value = 8000
$ws.Cells.Item(14, 13).Value = 7500
$ws.Cells.Item(14, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(14, 16).Value = 125
$ws.Cells.Item(14, 17).Value = 60

$ws.Cells.Item(15, 4).Value = 44321
$ws.Cells.Item(15, 10).Value = 500
$ws.Cells.Item(15, 11).Value = 7000
$ws.Cells.Item(15, 12).Value = 8000
$ws.Cells.Item(15, 13).Value = 7500
$ws.Cells.Item(15, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(15, 16).Value = 125
$ws.Cells.Item(15, 17).Value = 60

$ws.Cells.Item(16, 4).Value = 44517
$ws.Cells.Item(16, 10).Value = 500
$ws.Cells.Item(16, 11).Value = 5000
$ws.Cells.Item(16, 12).Value = 6000
$ws.Cells.Item(16, 13).Value = 5500
$ws.Cells.Item(16, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(16, 16).Value = 92
$ws.Cells.Item(16, 17).Value = 60

$ws.Cells.Item(17, 4).Value = 44265
$ws.Cells.Item(17, 10).Value = 500
$ws.Cells.Item(17, 11).Value = 6500
$ws.Cells.Item(17, 12).Value = 7000
$ws.Cells.Item(17, 13).Value = 6750
$ws.Cells.Item(17, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(17, 16).Value = 112
$ws.Cells.Item(17, 17).Value = 60

$ws.Cells.Item(18, 4).Value = 44328
$ws.Cells.Item(18, 10).Value = 500
$ws.Cells.Item(18, 11).Value = 7500
$ws.Cells.Item(18, 12).Value = 8000
$ws.Cells.Item(18, 13).Value = 7750
$ws.Cells.Item(18, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(18, 16).Value = 155
$ws.Cells.Item(18, 17).Value = 50

$ws.Cells.Item(19, 4).Value = 44216
$ws.Cells.Item(19, 10).Value = 1100
$ws.Cells.Item(19, 11).Value = 5500
$ws.Cells.Item(19, 12).Value = 6000
$ws.Cells.Item(19, 13).Value = 5750
$ws.Cells.Item(19, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(19, 16).Value = 96
$ws.Cells.Item(19, 17).Value = 60

$ws.Cells.Item(20, 4).Value = 44308
$ws.Cells.Item(20, 10).Value = 400
$ws.Cells.Item(20, 11).Value = 6000
$ws.Cells.Item(20, 12).Value = 7000
$ws.Cells.Item(20, 13).Value = 6500
$ws.Cells.Item(20, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(20, 16).Value = 108
$ws.Cells.Item(20, 17).Value = 60

$ws.Cells.Item(21, 4).Value = 44293
$ws.Cells.Item(21, 10).Value = 400
$ws.Cells.Item(21, 11).Value = 7000
$ws.Cells.Item(21, 12).Value = 8000
$ws.Cells.Item(21, 13).Value = 7500
$ws.Cells.Item(21, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(21, 16).Value = 125
$ws.Cells.Item(21, 17).Value = 60

$ws.Cells.Item(22, 4).Value = 44336
$ws.Cells.Item(22, 10).Value = 600
$ws.Cells.Item(22, 11).Value = 8500
$ws.Cells.Item(22, 12).Value = 9000
$ws.Cells.Item(22, 13).Value = 8750
$ws.Cells.Item(22, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(22, 16).Value = 175
$ws.Cells.Item(22, 17).Value = 50

$ws.Cells.Item(23, 4).Value = 44510
$ws.Cells.Item(23, 10).Value = 900
$ws.Cells.Item(23, 11).Value = 5000
$ws.Cells.Item(23, 12).Value = 6000
$ws.Cells.Item(23, 13).Value = 5500
$ws.Cells.Item(23, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(23, 16).Value = 92
$ws.Cells.Item(23, 17).Value = 60

$ws.Cells.Item(24, 4).Value = 44503
$ws.Cells.Item(24, 10).Value = 1100
$ws.Cells.Item(24, 11).Value = 6500
$ws.Cells.Item(24, 12).Value = 7000
$ws.Cells.Item(24, 13).Value = 6750
$ws.Cells.Item(24, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(24, 16).Value = 112
$ws.Cells.Item(24, 17).Value = 60

$ws.Cells.Item(25, 4).Value = 44371
$ws.Cells.Item(25, 10).Value = 300
$ws.Cells.Item(25, 11).Value = 8500
$ws.Cells.Item(25, 12).Value = 9000
$ws.Cells.Item(25, 13).Value = 8750
$ws.Cells.Item(25, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(25, 16).Value = 146
$ws.Cells.Item(25, 17).Value = 60

$ws.Cells.Item(26, 4).Value = 44384
$ws.Cells.Item(26, 10).Value = 300
$ws.Cells.Item(26, 11).Value = 7000
$ws.Cells.Item(26, 12).Value = 8000
$ws.Cells.Item(26, 13).Value = 7500
$ws.Cells.Item(26, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(26, 16).Value = 125
$ws.Cells.Item(26, 17).Value = 60

$ws.Cells.Item(27, 4).Value = 44315
$ws.Cells.Item(27, 10).Value = 500
$ws.Cells.Item(27, 11).Value = 7000
$ws.Cells.Item(27, 12).Value = 8000
$ws.Cells.Item(27, 13).Value = 7500
$ws.Cells.Item(27, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(27, 16).Value = 125
$ws.Cells.Item(27, 17).Value = 60

$ws.Cells.Item(28, 4).Value = 44244
$ws.Cells.Item(28, 10).Value = 500
$ws.Cells.Item(28, 11).Value = 5000
$ws.Cells.Item(28, 12).Value = 6000
$ws.Cells.Item(28, 13).Value = 5500
$ws.Cells.Item(28, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(28, 16).Value = 92
$ws.Cells.Item(28, 17).Value = 60

$ws.Cells.Item(29, 4).Value = 44377
$ws.Cells.Item(29, 10).Value = 400
$ws.Cells.Item(29, 11).Value = 7000
$ws.Cells.Item(29, 12).Value = 8000
$ws.Cells.Item(29, 13).Value = 7500
$ws.Cells.Item(29, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(29, 16).Value = 125
$ws.Cells.Item(29, 17).Value = 60

$ws.Cells.Item(30, 4).Value = 44300
$ws.Cells.Item(30, 10).Value = 400
$ws.Cells.Item(30, 11).Value = 6000
$ws.Cells.Item(30, 12).Value = 7000
$ws.Cells.Item(30, 13).Value = 6500
$ws.Cells.Item(30, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(30, 16).Value = 108
$ws.Cells.Item(30, 17).Value = 60

$ws.Cells.Item(31, 4).Value = 44497
$ws.Cells.Item(31, 10).Value = 800
$ws.Cells.Item(31, 11).Value = 7500
$ws.Cells.Item(31, 12).Value = 8000
$ws.Cells.Item(31, 13).Value = 7750
$ws.Cells.Item(31, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(31, 16).Value = 129
$ws.Cells.Item(31, 17).Value = 60

$ws.Cells.Item(32, 4).Value = 44525
$ws.Cells.Item(32, 10).Value = 400
$ws.Cells.Item(32, 11).Value = 5000
$ws.Cells.Item(32, 12).Value = 6000
$ws.Cells.Item(32, 13).Value = 5500
$ws.Cells.Item(32, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(32, 16).Value = 92
$ws.Cells.Item(32, 17).Value = 60

$ws.Cells.Item(33, 4).Value = 44286
$ws.Cells.Item(33, 10).Value = 600
$ws.Cells.Item(33, 11).Value = 7000
$ws.Cells.Item(33, 12).Value = 8000
$ws.Cells.Item(33, 13).Value = 7500
$ws.Cells.Item(33, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(33, 16).Value = 150
$ws.Cells.Item(33, 17).Value = 50

$ws.Cells.Item(34, 4).Value = 44307
$ws.Cells.Item(34, 10).Value = 700
$ws.Cells.Item(34, 11).Value = 6000
$ws.Cells.Item(34, 12).Value = 7000
$ws.Cells.Item(34, 13).Value = 6500
$ws.Cells.Item(34, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(34, 16).Value = 108
$ws.Cells.Item(34, 17).Value = 60

$ws.Cells.Item(35, 4).Value = 44181
$ws.Cells.Item(35, 10).Value = 900
$ws.Cells.Item(35, 11).Value = 4500
$ws.Cells.Item(35, 12).Value = 5000
$ws.Cells.Item(35, 13).Value = 4750
$ws.Cells.Item(35, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(35, 16).Value = 79
$ws.Cells.Item(35, 17).Value = 60

$ws.Cells.Item(36, 4).Value = 44335
$ws.Cells.Item(36, 10).Value = 500
$ws.Cells.Item(36, 11).Value = 7500
$ws.Cells.Item(36, 12).Value = 8000
$ws.Cells.Item(36, 13).Value = 7750
$ws.Cells.Item(36, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(36, 16).Value = 155
$ws.Cells.Item(36, 17).Value = 50

$ws.Cells.Item(37, 4).Value = 44314
$ws.Cells.Item(37, 10).Value = 1100
$ws.Cells.Item(37, 11).Value = 7000
$ws.Cells.Item(37, 12).Value = 8000
$ws.Cells.Item(37, 13).Value = 7500
$ws.Cells.Item(37, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(37, 16).Value = 125
$ws.Cells.Item(37, 17).Value = 60

$ws.Cells.Item(38, 4).Value = 44490
$ws.Cells.Item(38, 10).Value = 600
$ws.Cells.Item(38, 11).Value = 13000
$ws.Cells.Item(38, 12).Value = 15000
$ws.Cells.Item(38, 13).Value = 14000
$ws.Cells.Item(38, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(38, 16).Value = 233
$ws.Cells.Item(38, 17).Value = 60

